$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B,C -> C,D, preserving their
# bestFit widths); a second "stimFile" column is being added next to "posFile".
$ws.Columns.Item(2).Insert()

# Remove the old row 3 (content is being consolidated into a single data row
# that carries both the social and non-social CS+/CS- stimuli).
$ws.Rows.Item(3).Delete()

# Row 1 - headers
$ws.Range("A1").Value = "posFile"
$ws.Range("B1").Value = "stimFile"
$ws.Range("C1").Value = "cs_plus_s"
$ws.Range("D1").Value = "cs_minus_s"
$ws.Range("E1").Value = "cs_plus_ns"
$ws.Range("F1").Value = "cs_minus_ns"

# Row 2 - data
$ws.Range("A2").Value = "positions.xlsx"
$ws.Range("B2").Value = "stimuli.xlsx"
$ws.Range("C2").Value = "stimuli/social/031_y_m_n_a.jpg"
$ws.Range("D2").Value = "stimuli/social/016_y_m_n_b.jpg"
$ws.Range("E2").Value = "stimuli/non-social/016_y_m_n_b_scrambled.jpg"
$ws.Range("F2").Value = "stimuli/non-social/031_y_m_n_a_scrambled.jpg"

# New column B gets (approximately) the same width as column A (posFile /
# stimFile are both short filename-style headers).
$ws.Columns.Item(2).ColumnWidth = 15.6

# Selection as left after editing
$ws.Range("B5").Select()
